# Tower_塔表.xlsx — add two new tower rows (1019 "暗龙娘1", 1020 "暗龙娘3")
# at rows 25/26, tweak the frozen-pane view/selection, and widen columns N/O.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 25 — tower id 1019 "暗龙娘1"
# ---------------------------------------------------------------------------
$ws.Cells.Item(25,1).Value  = 1019
$ws.Cells.Item(25,2).Value  = "暗龙娘1"
$ws.Cells.Item(25,3).Value  = "暗龙娘1"
$ws.Cells.Item(25,4).Value  = 160
$ws.Cells.Item(25,5).Value  = 1
$ws.Cells.Item(25,6).Value  = "1|2"
$ws.Cells.Item(25,7).Value  = 326237
$ws.Cells.Item(25,8).Value  = "Tower_attackTags_1|Tower_attackTags_2|Tower_attackTags_3|Tower_attackTags_4"
$ws.Cells.Item(25,9).Value  = "attackDamage|attackTime|attackCount|findRange"
$ws.Cells.Item(25,10).Value = "160|250|320"
$ws.Cells.Item(25,11).Value = "160|250|320"
$ws.Cells.Item(25,12).Value = "3A9D803A4C74C27DA0A11FA53B742510|E718B09E4408CE5534779780E5365B64|E456238842ACC53D8C01EAABD11B256C"
$ws.Cells.Item(25,13).Value = "07F582E14688F6E9C9F0D2951D0902E1"
$ws.Cells.Item(25,14).Value = 150927
$ws.Cells.Item(25,15).Value = 1.3
$ws.Cells.Item(25,16).Value = 0
$ws.Cells.Item(25,17).Value = "1.3|1.3|1.1"
$ws.Cells.Item(25,18).Value = "8|8|8"
$ws.Cells.Item(25,19).Value = "0|0|0"
$ws.Cells.Item(25,20).Value = "1|1.25|1.25"
$ws.Cells.Item(25,21).Value = "50|65|80"
$ws.Cells.Item(25,22).Value = 146119
$ws.Cells.Item(25,23).Value = "FF9C6CFF|65AAFFFF|FF45FEFF"
$ws.Cells.Item(25,24).Value = 279672
$ws.Cells.Item(25,25).Value = 1
$ws.Cells.Item(25,26).Value = "1|1"

# ---------------------------------------------------------------------------
# 2. Row 26 — tower id 1020 "暗龙娘3"
# ---------------------------------------------------------------------------
$ws.Cells.Item(26,1).Value  = 1020
$ws.Cells.Item(26,2).Value  = "暗龙娘3"
$ws.Cells.Item(26,3).Value  = "暗龙娘3"
$ws.Cells.Item(26,4).Value  = 70
$ws.Cells.Item(26,5).Value  = 1
$ws.Cells.Item(26,6).Value  = "1|2"
$ws.Cells.Item(26,7).Value  = 326236
$ws.Cells.Item(26,8).Value  = "Tower_attackTags_1|Tower_attackTags_2|Tower_attackTags_3|Tower_attackTags_4"
$ws.Cells.Item(26,9).Value  = "attackDamage|attackTime|attackCount|findRange"
$ws.Cells.Item(26,10).Value = "70|120|180"
$ws.Cells.Item(26,11).Value = "70|120|180"
$ws.Cells.Item(26,12).Value = "CD0581B64EAB0563BF8537BCE23878AE|E718B09E4408CE5534779780E5365B64|E456238842ACC53D8C01EAABD11B256C"
$ws.Cells.Item(26,13).Value = "2E34840A438C9605A8FB7C978CAFE3D9"
$ws.Cells.Item(26,14).Value = 181026
$ws.Cells.Item(26,15).Value = 1
$ws.Cells.Item(26,16).Value = 0
$ws.Cells.Item(26,17).Value = "0.6|0.6|0.4"
$ws.Cells.Item(26,18).Value = "1|1|1"
$ws.Cells.Item(26,19).Value = "0|0|0"
$ws.Cells.Item(26,20).Value = "1.5|1.6|1.6"
$ws.Cells.Item(26,21).Value = "55|80|110"
$ws.Cells.Item(26,22).Value = 144198
$ws.Cells.Item(26,23).Value = "FF9C6CFF|65AAFFFF|FF45FEFF"
$ws.Cells.Item(26,24).Value = 285283
$ws.Cells.Item(26,25).Value = 0.4
$ws.Cells.Item(26,26).Value = "1|1"

# ---------------------------------------------------------------------------
# 3. Style touch-ups — reuse existing cell formats (via copy/PasteSpecial of
#    formats only) instead of assigning .Style directly, so no new cellXfs
#    entries are minted; this mirrors formats already used elsewhere on the
#    sheet for the same kind of column.
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

$ws.Cells.Item(5,6).Copy() | Out-Null               # F (attackTags list) style
$ws.Cells.Item(25,6).PasteSpecial($xlPasteFormats) | Out-Null

$ws.Cells.Item(5,8).Copy() | Out-Null                # H (attackTags header) style
$ws.Range("H25:H26").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Cells.Item(5,9).Copy() | Out-Null                # I (attack field header) style
$ws.Cells.Item(25,9).PasteSpecial($xlPasteFormats) | Out-Null

$ws.Cells.Item(17,10).Copy() | Out-Null              # J25 style
$ws.Cells.Item(25,10).PasteSpecial($xlPasteFormats) | Out-Null

$ws.Cells.Item(18,11).Copy() | Out-Null              # K (highlighted) style
$ws.Range("K25:K26").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Cells.Item(18,11).Copy() | Out-Null              # J26 style (same as K)
$ws.Cells.Item(26,10).PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Widen columns N (14) and O (15); column P (16) keeps its old width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(14).ColumnWidth = 18.61
$ws.Columns.Item(15).ColumnWidth = 18.08

# ---------------------------------------------------------------------------
# 5. Frozen-pane view: scroll to column L and move the active selection to
#    M25 (was S14) in the frozen top-right pane.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 12
$ws.Range("M25").Select() | Out-Null
